$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B (star_count) and C (n) values for rows 2-48 (A column index 0-46 stays unchanged)
$values = @(
    @(6, 58),
    @(7, 32),
    @(8, 27),
    @(9, 13),
    @(10, 12),
    @(12, 11),
    @(11, 10),
    @(13, 10),
    @(27, 5),
    @(14, 5),
    @(15, 5),
    @(31, 4),
    @(18, 4),
    @(19, 4),
    @(30, 3),
    @(25, 3),
    @(26, 3),
    @(32, 3),
    @(23, 3),
    @(22, 3),
    @(20, 3),
    @(17, 3),
    @(34, 2),
    @(41, 2),
    @(48, 2),
    @(33, 1),
    @(37, 1),
    @(130, 1),
    @(113, 1),
    @(87, 1),
    @(85, 1),
    @(84, 1),
    @(83, 1),
    @(74, 1),
    @(70, 1),
    @(68, 1),
    @(67, 1),
    @(58, 1),
    @(55, 1),
    @(51, 1),
    @(28, 1),
    @(46, 1),
    @(44, 1),
    @(43, 1),
    @(42, 1),
    @(16, 1),
    @(264, 1)
)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $values[$i][0]
    $ws.Cells.Item($row, 3).Value = $values[$i][1]
}

# Remove the now-unused rows 49-72 (data previously went up to row 72)
$lastOldRow = 72
$firstRemoveRow = $startRow + $values.Length
if ($lastOldRow -ge $firstRemoveRow) {
    $deleteRange = $ws.Range($ws.Cells.Item($firstRemoveRow, 1), $ws.Cells.Item($lastOldRow, 1))
    $deleteRange.EntireRow.Delete()
}
